# Re-save refresh of the CustomSliceSectorColorsPieChart example workbook:
#  - bump the Aspose.Cells evaluation-watermark copyright year (2014 -> 2016)
#    on the "Evaluation Warning" sheet
#  - drop the explicit PaperSize from both sheets' page setup (left at the
#    application default instead of being pinned to A4/Letter id 9)

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Sheet1")
$wsWarn = $wb.Worksheets.Item("Evaluation Warning")

# Update the evaluation watermark text (shared string used by A5 on the
# "Evaluation Warning" sheet) to reflect the newer copyright year.
$wsWarn.Range("A5").Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2016 Aspose Pty Ltd."

# Clear the explicit paper size on both sheets' page setup.
$wsData.PageSetup.PaperSize = $null
$wsWarn.PageSetup.PaperSize = $null
